$d = $word.ActiveDocument

$replacements = @(
    @("679÷7=", "979÷7="),
    @("668÷3=", "346÷9="),
    @("812÷8=", "293÷6="),
    @("215÷4=", "219÷6="),
    @("636÷5=", "615÷4="),
    @("758÷7=", "746÷4="),
    @("684÷4=", "664÷9="),
    @("529÷8=", "923÷7="),
    @("299÷7=", "615÷3="),
    @("933÷8=", "276÷4="),
    @("824÷3=", "403÷4="),
    @("851÷4=", "695÷3="),
    @("975÷7=", "637÷2="),
    @("146÷2=", "766÷5="),
    @("465÷5=", "153÷3="),
    @("499÷4=", "549÷6="),
    @("497÷9=", "651÷2="),
    @("892÷3=", "104÷2="),
    @("538÷4=", "496÷8="),
    @("619÷9=", "633÷7="),
    @("499÷3=", "594÷8="),
    @("353÷5=", "376÷8="),
    @("455÷6=", "901÷8="),
    @("502÷7=", "262÷8="),
    @("522÷6=", "264÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
